$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.237.77"
$ws.Range("E2").Value = "  -2.49%  "
$ws.Range("D3").Value = "3.069.93"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'523.31"
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("D6").Value = "'135.59"
$ws.Range("E6").Value = "  -5.52%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.068.19"
$ws.Range("E8").Value = "  -2.37%  "
$ws.Range("D9").Value = "'0.472"
$ws.Range("E9").Value = "  +4.52%  "
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  -3.70%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "'0.137"
$ws.Range("E13").Value = "  +1.65%  "
$ws.Range("D14").Value = "3.599.48"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "'25.07"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("E16").Value = "  -4.86%  "
$ws.Range("D17").Value = "57.237.81"
$ws.Range("E17").Value = "  -2.48%  "
$ws.Range("D18").Value = "3.060.20"
$ws.Range("E18").Value = "  -2.35%  "
$ws.Range("D19").Value = "'5.84"
$ws.Range("E19").Value = "  -4.98%  "
$ws.Range("E20").Value = "  -4.19%  "
$ws.Range("E21").Value = "  -2.71%  "
$ws.Range("D22").Value = "'347.44"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").Value = "'68.68"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").Value = "'0.165"
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("E28").Value = "  -10.56%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -6.15%  "
$ws.Range("E31").Value = "  -3.31%  "
$ws.Range("D32").Value = "'5.83"
$ws.Range("E32").Value = "  -10.33%  "
$ws.Range("D33").Value = "'20.97"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("D34").Value = "'158.73"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("E36").Value = "  -7.81%  "
$ws.Range("D37").Value = "'5.96"
$ws.Range("E37").Value = "  -5.11%  "
$ws.Range("D38").Value = "'25.48"
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("D41").Value = "'1.55"
$ws.Range("E41").Value = "  -7.04%  "
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("D43").Value = "'0.689"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").Value = "2.410.37"
$ws.Range("E44").Value = "  +4.38%  "
$ws.Range("D45").Value = "'36.68"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "3.108.61"
$ws.Range("E47").Value = "  -2.28%  "
$ws.Range("E48").Value = "  -2.78%  "
$ws.Range("E49").Value = "  -2.36%  "
$ws.Range("E50").Value = "  -7.65%  "
$ws.Range("D51").Value = "'19.27"
$ws.Range("E51").Value = "  -7.26%  "
